$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) "...flow-routing algorithms from the pygeoprocessing module (cite)..."
#    becomes
#    "...flow-routing algorithms from the whitebox tools open-source
#     geoprocessing library (cite)..."
#    The new wording is split across several runs, mirroring how Word
#    segments text when it is retyped in place (the existing spell-check
#    proofErr wrapper around the package name is preserved).
# ---------------------------------------------------------------------------

# 1a) Split "...flow-routing algorithms from the " into
#     "...flow-routing algorithms from " + "the "
$fr = $d.Content
$fr.Find.Execute("flow-routing algorithms from the ") | Out-Null
$fr.Text = "flow-routing algorithms from "
$fr.Collapse(0)
$insertStart = $fr.Start
$fr.InsertAfter("the ")
$newRun = $d.Range($insertStart, $insertStart + 4)
# Toggling the font size and restoring it forces Word to materialize this
# span as its own run instead of merging it back into its neighbour.
$newRun.Font.Size = 99
$newRun.Font.Size = 12

# 1b) pygeoprocessing -> whitebox (keeps its spellStart/spellEnd wrapper)
$d.Content.Find.Execute("pygeoprocessing", $true, $false, $false, $false, `
    $false, $true, 1, $false, "whitebox", 2) | Out-Null

# 1c) Split " module (cite). K factors were extracted " into
#     " tools open-source geoprocessing library" + " (cite). K factors were extracted "
$fr2 = $d.Content
$fr2.Find.Execute(" module (cite). K factors were extracted ") | Out-Null
$fr2.Text = " tools open-source geoprocessing library (cite). K factors were extracted "
$splitStart = $fr2.Start
$firstPart = " tools open-source geoprocessing library"
$firstRun = $d.Range($splitStart, $splitStart + $firstPart.Length)
$firstRun.Font.Size = 99
$firstRun.Font.Size = 12

# ---------------------------------------------------------------------------
# 2) Move the "_GoBack" last-edit bookmark from the end of the document to
#    right after the "Modelled" entry in the "Soil test Aluminum" row of the
#    parameter table -- i.e. where the text above was actually edited.
# ---------------------------------------------------------------------------
$table = $d.Tables.Item(1)
$cell = $table.Cell(7, 2)
$cellRange = $cell.Range
$lastCharPos = $cellRange.End - 1
$target = $d.Range($lastCharPos, $lastCharPos)
$target.InsertAfter("X")
$d.Bookmarks.Add("_GoBack", $target) | Out-Null
$target.Text = ""
